$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M6")

# Add new rows of data (continuing the existing log table), per
# "Aufgabe 1.7 Viele sachen aus M5 hinzugefügt, die leider fehlten."
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = "css probleme"

$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 20
$ws.Range("E11").Value = "Flow war da"

# Move/update the active cell selection to the next empty row, as in the
# saved workbook state.
$ws.Range("B12").Select()
